$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 175
$ws.Range("I12").Value = 175
$ws.Range("K12").Value = 175
$ws.Range("M12").Value = -5
$ws.Range("H17").Value = 500
$ws.Range("J17").Value = 566.6667
$ws.Range("L17").Value = 1700.0001
$ws.Range("N17").Value = -2036.0001
$ws.Range("H86").Value = 9857.143
$ws.Range("I86").Value = 9000.25
$ws.Range("J86").Value = 10999.667
$ws.Range("K86").Value = 9000.25
$ws.Range("L86").Value = 10999.667
$ws.Range("M86").Value = -7877.25
$ws.Range("N86").Value = -13245.667
$ws.Range("H88").Value = 2579.4167
$ws.Range("I88").Value = 1761
$ws.Range("J88").Value = 3164
$ws.Range("K88").Value = 1761
$ws.Range("L88").Value = 3164
$ws.Range("M88").Value = -1355
$ws.Range("N88").Value = -3976
$ws.Range("H89").Value = 9857.143
$ws.Range("I89").Value = 9000.25
$ws.Range("J89").Value = 10999.667
$ws.Range("K89").Value = 45001.25
$ws.Range("L89").Value = 54998.335
$ws.Range("M89").Value = -39385.25
$ws.Range("N89").Value = -66230.33499999999
$ws.Range("H91").Value = 2579.4167
$ws.Range("I91").Value = 1761
$ws.Range("J91").Value = 3164
$ws.Range("K91").Value = 1761
$ws.Range("L91").Value = 3164
$ws.Range("M91").Value = -357
$ws.Range("N91").Value = -5972
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents() | Out-Null
$ws.Range("H129").Value = 1098.5
$ws.Range("I129").Value = 1098.5
$ws.Range("K129").Value = 3295.5
$ws.Range("M129").Value = 1704.5
$ws.Range("H132").Value = 1637.5385
$ws.Range("I132").Value = 1637.5385
$ws.Range("K132").Value = 4912.6155
$ws.Range("M132").Value = -2382.6155
$ws.Range("H138").Value = 1549.1765
$ws.Range("I138").Value = 1301.4
$ws.Range("K138").Value = 3904.2
$ws.Range("M138").Value = 1235.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 1225.75
$ws.Range("I6").Value = 1225.75
$ws.Range("K6").Value = 1225.75
$ws.Range("M6").Value = -1052.75
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents() | Out-Null
$ws.Range("H32").Value = 3840.2104
$ws.Range("I32").Value = 3840.2104
$ws.Range("K32").Value = 3840.2104
$ws.Range("M32").Value = -3553.2104
$ws.Range("H45").Value = 2627.7144
$ws.Range("I45").Value = 1199.5
$ws.Range("J45").Value = 3199
$ws.Range("K45").Value = 1199.5
$ws.Range("L45").Value = 3199
$ws.Range("M45").Value = -822.5
$ws.Range("N45").Value = -3953
$ws.Range("H97").Value = 1487.7222
$ws.Range("I97").Value = 1221.7273
$ws.Range("J97").Value = 1905.7142
$ws.Range("K97").Value = 1221.7273
$ws.Range("L97").Value = 1905.7142
$ws.Range("M97").Value = -725.7273
$ws.Range("N97").Value = -2897.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1082.5
$ws.Range("I20").Value = 899
$ws.Range("J20").Value = 2000
$ws.Range("K20").Value = 899
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = -652
$ws.Range("N20").Value = -2494
$ws.Range("H86").Value = 2109.5
$ws.Range("I86").Value = 704.6
$ws.Range("K86").Value = 704.6
$ws.Range("M86").Value = 418.4
$ws.Range("H89").Value = 2109.5
$ws.Range("I89").Value = 704.6
$ws.Range("K89").Value = 3523
$ws.Range("M89").Value = 2093
$ws.Range("H94").Value = 2605.4
$ws.Range("I94").Value = 2605.4
$ws.Range("K94").Value = 2605.4
$ws.Range("M94").Value = -2154.4
$ws.Range("H105").Value = 702.4
$ws.Range("I105").Value = 706.75
$ws.Range("K105").Value = 706.75
$ws.Range("M105").Value = 1040.25
$ws.Range("H107").Value = 784.8
$ws.Range("I107").Value = 991.3333
$ws.Range("J107").Value = 475
$ws.Range("K107").Value = 991.3333
$ws.Range("L107").Value = 475
$ws.Range("M107").Value = 928.6667
$ws.Range("N107").Value = -4315

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1499.5
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 1499
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 4497
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -9557
$ws.Range("H134").Value = 1650
$ws.Range("J134").Value = 1650
$ws.Range("L134").Value = 4950
$ws.Range("N134").Value = -10020

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3722.3076
$ws.Range("I68").Value = 3578.2222
$ws.Range("J68").Value = 4046.5
$ws.Range("K68").Value = 10734.6666
$ws.Range("L68").Value = 12139.5
$ws.Range("M68").Value = -9923.6666
$ws.Range("N68").Value = -13761.5
$ws.Range("H71").Value = 3722.3076
$ws.Range("I71").Value = 3578.2222
$ws.Range("J71").Value = 4046.5
$ws.Range("K71").Value = 32203.9998
$ws.Range("L71").Value = 36418.5
$ws.Range("M71").Value = -28147.9998
$ws.Range("N71").Value = -44530.5
$ws.Range("H131").Value = 986.7742
$ws.Range("I131").Value = 890
$ws.Range("K131").Value = 2670
$ws.Range("M131").Value = 2370

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6397.143
$ws.Range("I80").Value = 3295.3333
$ws.Range("J80").Value = 8723.5
$ws.Range("K80").Value = 3295.3333
$ws.Range("L80").Value = 8723.5
$ws.Range("M80").Value = -2297.3333
$ws.Range("N80").Value = -10719.5
$ws.Range("H83").Value = 6397.143
$ws.Range("I83").Value = 3295.3333
$ws.Range("J83").Value = 8723.5
$ws.Range("K83").Value = 16476.6665
$ws.Range("L83").Value = 43617.5
$ws.Range("M83").Value = -11484.6665
$ws.Range("N83").Value = -53601.5
$ws.Range("H113").Value = 5086.4165
$ws.Range("I113").Value = 4582
$ws.Range("K113").Value = 4582
$ws.Range("M113").Value = -2412

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1269.6154
$ws.Range("J82").Value = 851.5
$ws.Range("L82").Value = 851.5
$ws.Range("N82").Value = -1573.5
$ws.Range("H85").Value = 1269.6154
$ws.Range("J85").Value = 851.5
$ws.Range("L85").Value = 851.5
$ws.Range("N85").Value = -3347.5
$ws.Range("H132").Value = 2774.8333
$ws.Range("I132").Value = 2774.8333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8324.499899999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5794.499899999999
$ws.Range("N132").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 66199.60000000001
$ws.Range("J119").Value = 66199.60000000001
$ws.Range("L119").Value = 66199.60000000001
$ws.Range("N119").Value = -75875.60000000001
$ws.Range("H122").Value = 1605.8636
$ws.Range("I122").Value = 1697.2106
$ws.Range("J122").Value = 1027.3334
$ws.Range("K122").Value = 5091.6318
$ws.Range("L122").Value = 3082.0002
$ws.Range("M122").Value = -2641.6318
$ws.Range("N122").Value = -7982.0002
